$wb = $excel.ActiveWorkbook

# Remember the sheet that is currently active so we can restore selection/
# activation state once we're done (adding/activating a new sheet would
# otherwise shift the workbook's active tab).
$origActive = $wb.ActiveSheet

# The style used on the header row of the other "summary" sheets (bold,
# centered, bordered) so we can copy it onto the new sheet's header.
$styleSource = $wb.Worksheets.Item("Proportion Summary").Range("A1")

# Add the new worksheet after the last existing sheet (i.e. at the end of
# the tab strip).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Num Quotes Proportions"

# Header row
$ws.Range("A1").Value = "year"
$ws.Range("B1").Value = "proportion_num_quotes_to_num_articles"

# Match the header formatting used elsewhere in the workbook.
$styleSource.Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)

# Data rows: year -> proportion of number of quotes to number of articles
$data = @(
    @(2019, 7.25),
    @(2020, 5.142857142857143),
    @(2021, 1.4),
    @(2022, 4.833333333333333),
    @(2023, 4),
    @(2024, 1.333333333333333)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row++
}

# Restore the original active sheet/selection.
$origActive.Activate()
